# Add data for 2022-09-30 (rolling window now "through September 22")
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet tab to reflect the new "as of" date.
$ws.Name = "Through 2022-09-22"

# Update the header label for the current (partial) month column (B1),
# which is backed by the same shared string as the sheet name context.
$ws.Range("B1").Value = "September 2022 (through September 22)"

# Per-neighborhood carjacking counts for the current partial month (column B)
# and other affected months, incremented by the new day's incidents.
$ws.Range("B2").Value  = 6    # Austin
$ws.Range("K2").Value  = 9    # Austin
$ws.Range("T2").Value  = 10   # Austin
$ws.Range("AU2").Value = 5    # Austin

$ws.Range("B6").Value  = 5    # Englewood

$ws.Range("K7").Value  = 4    # Little Italy, UIC
$ws.Range("AC7").Value = 1    # Little Italy, UIC (new data point)

$ws.Range("K10").Value = 8    # North Lawndale

$ws.Range("B11").Value = 2    # West Pullman

$ws.Range("B12").Value  = 3   # Grand Crossing
$ws.Range("AL12").Value = 4   # Grand Crossing

$ws.Range("T14").Value = 4    # Roseland

$ws.Range("K15").Value = 6    # Chatham

$ws.Range("AC23").Value = 3   # South Chicago

$ws.Range("K33").Value = 4    # United Center

$ws.Range("B34").Value = 1    # Uptown (new data point)

$ws.Range("AC38").Value = 3   # Wicker Park

$ws.Range("AL45").Value = 1   # Clearing (new data point)

$ws.Range("AC50").Value = 1   # Grand Boulevard (new data point)

$ws.Range("K60").Value = 2    # Armour Square

$ws.Range("AC64").Value = 2   # Douglas

$ws.Range("AU74").Value = 2   # Irving Park
